# Update cryptos list with refreshed prices and 1h volume changes
# (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.501.24"
Set-TextValue $ws.Range("E2") "  +1.69%  "
Set-TextValue $ws.Range("D3") "3.941.38"
Set-TextValue $ws.Range("E3") "  +0.24%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.22%  "
Set-TextValue $ws.Range("D5") "505.37"
Set-TextValue $ws.Range("E5") "  +3.72%  "
Set-TextValue $ws.Range("D6") "147.72"
Set-TextValue $ws.Range("E6") "  -0.53%  "
Set-TextValue $ws.Range("D7") "0.625"
Set-TextValue $ws.Range("E7") "  -0.50%  "
Set-TextValue $ws.Range("E8") "  +0.02%  "
Set-TextValue $ws.Range("D9") "0.735"
Set-TextValue $ws.Range("E9") "  -0.19%  "
Set-TextValue $ws.Range("D10") "0.176"
Set-TextValue $ws.Range("E10") "  +3.95%  "
Set-TextValue $ws.Range("E11") "  -0.91%  "
Set-TextValue $ws.Range("D12") "43.51"
Set-TextValue $ws.Range("E12") "  +1.01%  "
Set-TextValue $ws.Range("D13") "10.51"
Set-TextValue $ws.Range("E13") "  -1.79%  "
Set-TextValue $ws.Range("D14") "4.573.83"
Set-TextValue $ws.Range("E14") "  +0.25%  "
Set-TextValue $ws.Range("D15") "3.941.32"
Set-TextValue $ws.Range("E15") "  +0.56%  "
Set-TextValue $ws.Range("D16") "14.23"
Set-TextValue $ws.Range("E16") "  -2.23%  "
Set-TextValue $ws.Range("E17") "  -0.27%  "
Set-TextValue $ws.Range("E18") "  +5.60%  "
Set-TextValue $ws.Range("D19") "20.00"
Set-TextValue $ws.Range("E19") "  -0.07%  "
Set-TextValue $ws.Range("D20") "69.478.80"
Set-TextValue $ws.Range("D21") "436.50"
Set-TextValue $ws.Range("E21") "  -1.45%  "
Set-TextValue $ws.Range("D22") "3.44"
Set-TextValue $ws.Range("E22") "  -2.25%  "
Set-TextValue $ws.Range("D23") "14.72"
Set-TextValue $ws.Range("E23") "  -2.48%  "
Set-TextValue $ws.Range("D24") "89.00"
Set-TextValue $ws.Range("E24") "  +0.58%  "
Set-TextValue $ws.Range("D25") "11.97"
Set-TextValue $ws.Range("E25") "  +5.08%  "
Set-TextValue $ws.Range("E26") "  +6.35%  "
Set-TextValue $ws.Range("D27") "11.23"
Set-TextValue $ws.Range("E27") "  -1.92%  "
Set-TextValue $ws.Range("D28") "37.17"
Set-TextValue $ws.Range("E28") "  -4.70%  "
Set-TextValue $ws.Range("D29") "5.66"
Set-TextValue $ws.Range("E29") "  -2.96%  "
Set-TextValue $ws.Range("D30") "706.71"
Set-TextValue $ws.Range("E30") "  -2.15%  "
Set-TextValue $ws.Range("E31") "  -2.02%  "
Set-TextValue $ws.Range("E32") "  -1.56%  "
Set-TextValue $ws.Range("E33") "  -0.77%  "
Set-TextValue $ws.Range("D34") "64.75"
Set-TextValue $ws.Range("E34") "  +5.64%  "
Set-TextValue $ws.Range("D35") "0.452"
Set-TextValue $ws.Range("E35") "  +12.36%  "
Set-TextValue $ws.Range("D36") "0.0₃0889"
Set-TextValue $ws.Range("E36") "  -0.92%  "
Set-TextValue $ws.Range("D37") "6.07"
Set-TextValue $ws.Range("E37") "  -3.43%  "
Set-TextValue $ws.Range("D38") "40.94"
Set-TextValue $ws.Range("E38") "  -3.48%  "
Set-TextValue $ws.Range("E39") "  +1.67%  "
Set-TextValue $ws.Range("D40") "0.999"
Set-TextValue $ws.Range("E40") "  -0.02%  "
Set-TextValue $ws.Range("E41") "  -0.03%  "
Set-TextValue $ws.Range("E42") "  +1.62%  "
Set-TextValue $ws.Range("D43") "2.89"
Set-TextValue $ws.Range("E43") "  -5.34%  "
Set-TextValue $ws.Range("B44") "WEMIXToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("E44") "  +4.75%  "
Set-TextValue $ws.Range("B45") "ThetaToken"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D45") "3.08"
Set-TextValue $ws.Range("E45") "  -4.74%  "
Set-TextValue $ws.Range("E46") "  +1.07%  "
Set-TextValue $ws.Range("D47") "3.37"
Set-TextValue $ws.Range("E47") "  +3.68%  "
Set-TextValue $ws.Range("E48") "  +5.76%  "
Set-TextValue $ws.Range("D49") "3.40"
Set-TextValue $ws.Range("E49") "  -0.69%  "
Set-TextValue $ws.Range("D50") "0.0₆0350"
Set-TextValue $ws.Range("E50") "  -6.17%  "
Set-TextValue $ws.Range("D51") "2.10"
Set-TextValue $ws.Range("E51") "  -2.64%  "
